$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Shift the two existing trailing rows (204, 205) down to (206, 207)
#    by inserting two fresh rows at position 204. This preserves the
#    original cell formatting (bold/bordered id column, date format)
#    of the rows being pushed down.
# ------------------------------------------------------------------
$ws.Range("A204:A205").EntireRow.Insert()

# ------------------------------------------------------------------
# 2) Clone formatting (bold+border id cell, date number format) from
#    row 203 onto the two brand-new rows 204:205 so new cells line up
#    with the existing style table (s=1 / s=2) instead of Excel
#    minting new style entries.
# ------------------------------------------------------------------
$ws.Range("A203:G203").Copy()
$ws.Range("A204:G205").PasteSpecial(-4122)
$ws.Range("K203:AC203").Copy()
$ws.Range("K204:AC205").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Populate the two newly inserted rows (204, 205) with the two
#    earlier-dated fixtures that now slot in before the previous
#    last two rows.
# ------------------------------------------------------------------
# Row 204
$ws.Cells.Item(204,1).Value = 202
$ws.Cells.Item(204,2).Value = 6788918
$ws.Cells.Item(204,3).Value = "Croatia HNL"
$ws.Cells.Item(204,4).Value = "Croatia HNL"
$ws.Cells.Item(204,5).Value = 45332.45833333334
$ws.Cells.Item(204,6).Value = "Istra 1961"
$ws.Cells.Item(204,7).Value = "NK Osijek"
$ws.Cells.Item(204,8).Value = 1
$ws.Cells.Item(204,9).Value = 0
$ws.Cells.Item(204,10).Value = "H"
$ws.Cells.Item(204,11).Value = 3.25
$ws.Cells.Item(204,12).Value = 3.4
$ws.Cells.Item(204,13).Value = 2
$ws.Cells.Item(204,14).Value = 3.1
$ws.Cells.Item(204,15).Value = 2.9
$ws.Cells.Item(204,16).Value = 2.25
$ws.Cells.Item(204,17).Value = 0.25
$ws.Cells.Item(204,18).Value = 1.825
$ws.Cells.Item(204,19).Value = 2.025
$ws.Cells.Item(204,20).Value = 2
$ws.Cells.Item(204,21).Value = 1.8
$ws.Cells.Item(204,22).Value = 2.05
$ws.Cells.Item(204,23).Value = 2.1
$ws.Cells.Item(204,24).Value = -1
$ws.Cells.Item(204,25).Value = -1
$ws.Cells.Item(204,26).Value = 0.825
$ws.Cells.Item(204,27).Value = -1
$ws.Cells.Item(204,28).Value = -1
$ws.Cells.Item(204,29).Value = 1.05

# Row 205
$ws.Cells.Item(205,1).Value = 203
$ws.Cells.Item(205,2).Value = 6787891
$ws.Cells.Item(205,3).Value = "Croatia HNL"
$ws.Cells.Item(205,4).Value = "Croatia HNL"
$ws.Cells.Item(205,5).Value = 45332.54861111111
$ws.Cells.Item(205,6).Value = "NK Lokomotiva Zagreb"
$ws.Cells.Item(205,7).Value = "Dinamo Zagreb"
$ws.Cells.Item(205,8).Value = 2
$ws.Cells.Item(205,9).Value = 2
$ws.Cells.Item(205,10).Value = "D"
$ws.Cells.Item(205,11).Value = 5.5
$ws.Cells.Item(205,12).Value = 3.75
$ws.Cells.Item(205,13).Value = 1.615
$ws.Cells.Item(205,14).Value = 5.5
$ws.Cells.Item(205,15).Value = 3.75
$ws.Cells.Item(205,16).Value = 1.615
$ws.Cells.Item(205,17).Value = 0.75
$ws.Cells.Item(205,18).Value = 2.05
$ws.Cells.Item(205,19).Value = 1.8
$ws.Cells.Item(205,20).Value = 2.5
$ws.Cells.Item(205,21).Value = 1.975
$ws.Cells.Item(205,22).Value = 1.875
$ws.Cells.Item(205,23).Value = -1
$ws.Cells.Item(205,24).Value = 2.75
$ws.Cells.Item(205,25).Value = -1
$ws.Cells.Item(205,26).Value = 1.05
$ws.Cells.Item(205,27).Value = -1
$ws.Cells.Item(205,28).Value = 0.9750000000000001
$ws.Cells.Item(205,29).Value = -1

# ------------------------------------------------------------------
# 4) The rows that used to be 204/205 are now 206/207: refresh their
#    id numbers and fill in the now-known match result (FTHG/FTAG/FTR)
#    plus the recalculated odds columns.
# ------------------------------------------------------------------
# Row 206
$ws.Cells.Item(206,1).Value = 204
$ws.Cells.Item(206,2).Value = 6788917
$ws.Cells.Item(206,3).Value = "Croatia HNL"
$ws.Cells.Item(206,4).Value = "Croatia HNL"
$ws.Cells.Item(206,5).Value = 45333.45833333334
$ws.Cells.Item(206,6).Value = "Hajduk Split"
$ws.Cells.Item(206,7).Value = "Slaven Belupo"
$ws.Cells.Item(206,8).Value = 4
$ws.Cells.Item(206,9).Value = 0
$ws.Cells.Item(206,10).Value = "H"
$ws.Cells.Item(206,11).Value = 1.25
$ws.Cells.Item(206,12).Value = 5.75
$ws.Cells.Item(206,13).Value = 9
$ws.Cells.Item(206,14).Value = 1.285
$ws.Cells.Item(206,15).Value = 5.5
$ws.Cells.Item(206,16).Value = 8.5
$ws.Cells.Item(206,17).Value = -1.5
$ws.Cells.Item(206,18).Value = 1.95
$ws.Cells.Item(206,19).Value = 1.9
$ws.Cells.Item(206,20).Value = 2.75
$ws.Cells.Item(206,21).Value = 2.025
$ws.Cells.Item(206,22).Value = 1.825
$ws.Cells.Item(206,23).Value = 0.2849999999999999
$ws.Cells.Item(206,24).Value = -1
$ws.Cells.Item(206,25).Value = -1
$ws.Cells.Item(206,26).Value = 0.95
$ws.Cells.Item(206,27).Value = -1
$ws.Cells.Item(206,28).Value = 1.025
$ws.Cells.Item(206,29).Value = -1

# Row 207
$ws.Cells.Item(207,1).Value = 205
$ws.Cells.Item(207,2).Value = 6788919
$ws.Cells.Item(207,3).Value = "Croatia HNL"
$ws.Cells.Item(207,4).Value = "Croatia HNL"
$ws.Cells.Item(207,5).Value = 45333.54861111111
$ws.Cells.Item(207,6).Value = "HNK Gorica"
$ws.Cells.Item(207,7).Value = "HNK Rijeka"
$ws.Cells.Item(207,8).Value = 0
$ws.Cells.Item(207,9).Value = 2
$ws.Cells.Item(207,10).Value = "A"
$ws.Cells.Item(207,11).Value = 4.75
$ws.Cells.Item(207,12).Value = 3.5
$ws.Cells.Item(207,13).Value = 1.727
$ws.Cells.Item(207,14).Value = 6.5
$ws.Cells.Item(207,15).Value = 3.8
$ws.Cells.Item(207,16).Value = 1.533
$ws.Cells.Item(207,17).Value = 1
$ws.Cells.Item(207,18).Value = 1.925
$ws.Cells.Item(207,19).Value = 1.925
$ws.Cells.Item(207,20).Value = 2.5
$ws.Cells.Item(207,21).Value = 1.975
$ws.Cells.Item(207,22).Value = 1.875
$ws.Cells.Item(207,23).Value = -1
$ws.Cells.Item(207,24).Value = -1
$ws.Cells.Item(207,25).Value = 0.5329999999999999
$ws.Cells.Item(207,26).Value = -1
$ws.Cells.Item(207,27).Value = 0.925
$ws.Cells.Item(207,28).Value = -1
$ws.Cells.Item(207,29).Value = 0.875

# ------------------------------------------------------------------
# 5) Append five brand-new upcoming fixtures as rows 208-212 (no
#    result columns yet, so only A:G and K:AA are populated).
# ------------------------------------------------------------------
$ws.Range("A203:G203").Copy()
$ws.Range("A208:G212").PasteSpecial(-4122)
$ws.Range("K203:AA203").Copy()
$ws.Range("K208:AA212").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 208
$ws.Cells.Item(208,1).Value = 206
$ws.Cells.Item(208,2).Value = 6788922
$ws.Cells.Item(208,3).Value = "Croatia HNL"
$ws.Cells.Item(208,4).Value = "Croatia HNL"
$ws.Cells.Item(208,5).Value = 45338.54166666666
$ws.Cells.Item(208,6).Value = "Slaven Belupo"
$ws.Cells.Item(208,7).Value = "Istra 1961"
$ws.Cells.Item(208,11).Value = 2
$ws.Cells.Item(208,12).Value = 3.2
$ws.Cells.Item(208,13).Value = 3.4
$ws.Cells.Item(208,14).Value = 2.3
$ws.Cells.Item(208,15).Value = 3
$ws.Cells.Item(208,16).Value = 3
$ws.Cells.Item(208,17).Value = -0.25
$ws.Cells.Item(208,18).Value = 2.05
$ws.Cells.Item(208,19).Value = 1.8
$ws.Cells.Item(208,20).Value = 2.25
$ws.Cells.Item(208,21).Value = 1.975
$ws.Cells.Item(208,22).Value = 1.875
$ws.Cells.Item(208,23).Value = 0
$ws.Cells.Item(208,24).Value = 0
$ws.Cells.Item(208,25).Value = 0
$ws.Cells.Item(208,26).Value = 0
$ws.Cells.Item(208,27).Value = 0

# Row 209
$ws.Cells.Item(209,1).Value = 207
$ws.Cells.Item(209,2).Value = 6787892
$ws.Cells.Item(209,3).Value = "Croatia HNL"
$ws.Cells.Item(209,4).Value = "Croatia HNL"
$ws.Cells.Item(209,5).Value = 45339.45833333334
$ws.Cells.Item(209,6).Value = "HNK Rijeka"
$ws.Cells.Item(209,7).Value = "NK Lokomotiva Zagreb"
$ws.Cells.Item(209,11).Value = 1.363
$ws.Cells.Item(209,12).Value = 4.5
$ws.Cells.Item(209,13).Value = 7
$ws.Cells.Item(209,14).Value = 1.4
$ws.Cells.Item(209,15).Value = 4.5
$ws.Cells.Item(209,16).Value = 6.5
$ws.Cells.Item(209,17).Value = -1.25
$ws.Cells.Item(209,18).Value = 1.975
$ws.Cells.Item(209,19).Value = 1.875
$ws.Cells.Item(209,20).Value = 2.75
$ws.Cells.Item(209,21).Value = 1.975
$ws.Cells.Item(209,22).Value = 1.875
$ws.Cells.Item(209,23).Value = 0
$ws.Cells.Item(209,24).Value = 0
$ws.Cells.Item(209,25).Value = 0
$ws.Cells.Item(209,26).Value = 0
$ws.Cells.Item(209,27).Value = 0

# Row 210
$ws.Cells.Item(210,1).Value = 208
$ws.Cells.Item(210,2).Value = 6769301
$ws.Cells.Item(210,3).Value = "Croatia HNL"
$ws.Cells.Item(210,4).Value = "Croatia HNL"
$ws.Cells.Item(210,5).Value = 45339.54861111111
$ws.Cells.Item(210,6).Value = "NK Rudes"
$ws.Cells.Item(210,7).Value = "Hajduk Split"
$ws.Cells.Item(210,11).Value = 7.5
$ws.Cells.Item(210,12).Value = 4.75
$ws.Cells.Item(210,13).Value = 1.333
$ws.Cells.Item(210,14).Value = 10
$ws.Cells.Item(210,15).Value = 5.5
$ws.Cells.Item(210,16).Value = 1.25
$ws.Cells.Item(210,17).Value = 1.75
$ws.Cells.Item(210,18).Value = 1.825
$ws.Cells.Item(210,19).Value = 2.025
$ws.Cells.Item(210,20).Value = 2.75
$ws.Cells.Item(210,21).Value = 1.925
$ws.Cells.Item(210,22).Value = 1.925
$ws.Cells.Item(210,23).Value = 0
$ws.Cells.Item(210,24).Value = 0
$ws.Cells.Item(210,25).Value = 0
$ws.Cells.Item(210,26).Value = 0
$ws.Cells.Item(210,27).Value = 0

# Row 211
$ws.Cells.Item(211,1).Value = 209
$ws.Cells.Item(211,2).Value = 6788921
$ws.Cells.Item(211,3).Value = "Croatia HNL"
$ws.Cells.Item(211,4).Value = "Croatia HNL"
$ws.Cells.Item(211,5).Value = 45340.45833333334
$ws.Cells.Item(211,6).Value = "NK Osijek"
$ws.Cells.Item(211,7).Value = "HNK Gorica"
$ws.Cells.Item(211,11).Value = 1.75
$ws.Cells.Item(211,12).Value = 3.5
$ws.Cells.Item(211,13).Value = 4.2
$ws.Cells.Item(211,14).Value = 1.75
$ws.Cells.Item(211,15).Value = 3.4
$ws.Cells.Item(211,16).Value = 4.333
$ws.Cells.Item(211,17).Value = -0.5
$ws.Cells.Item(211,18).Value = 1.825
$ws.Cells.Item(211,19).Value = 2.025
$ws.Cells.Item(211,20).Value = 2.5
$ws.Cells.Item(211,21).Value = 2
$ws.Cells.Item(211,22).Value = 1.85
$ws.Cells.Item(211,23).Value = 0
$ws.Cells.Item(211,24).Value = 0
$ws.Cells.Item(211,25).Value = 0
$ws.Cells.Item(211,26).Value = 0
$ws.Cells.Item(211,27).Value = 0

# Row 212
$ws.Cells.Item(212,1).Value = 210
$ws.Cells.Item(212,2).Value = 6788920
$ws.Cells.Item(212,3).Value = "Croatia HNL"
$ws.Cells.Item(212,4).Value = "Croatia HNL"
$ws.Cells.Item(212,5).Value = 45340.54861111111
$ws.Cells.Item(212,6).Value = "Dinamo Zagreb"
$ws.Cells.Item(212,7).Value = "NK Varazdin"
$ws.Cells.Item(212,11).Value = 1.166
$ws.Cells.Item(212,12).Value = 6
$ws.Cells.Item(212,13).Value = 15
$ws.Cells.Item(212,14).Value = 1.2
$ws.Cells.Item(212,15).Value = 5.75
$ws.Cells.Item(212,16).Value = 11
$ws.Cells.Item(212,17).Value = -1.75
$ws.Cells.Item(212,18).Value = 1.875
$ws.Cells.Item(212,19).Value = 1.975
$ws.Cells.Item(212,20).Value = 3
$ws.Cells.Item(212,21).Value = 1.975
$ws.Cells.Item(212,22).Value = 1.875
$ws.Cells.Item(212,23).Value = 0
$ws.Cells.Item(212,24).Value = 0
$ws.Cells.Item(212,25).Value = 0
$ws.Cells.Item(212,26).Value = 0
$ws.Cells.Item(212,27).Value = 0

# ------------------------------------------------------------------
# 6) Make sure the sheet dimension reflects the new extent.
# ------------------------------------------------------------------
$ws.Range("A1:AC212").Select()
$ws.Range("A1").Select()
